$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VO IDs range")

# Assign new ID for new VO terms: VO:0010458 -> VO:0010461
$ws.Range("A13").Value = "VO:0010461"

# RxNorm Term starting from: VO:0021167 -> VO:0021180
$ws.Range("A17").Value = "VO:0021180"

# Reflect the editor's final cursor position on the sheet.
[void]$ws.Range("A19").Select()
